$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D): values that Excel would otherwise auto-parse as a
# number (single decimal point, e.g. "562.16"). Force them to remain plain
# text by applying a Text number format before assigning the value, then
# restoring the original "Normal" style so no visible formatting changes.
$forceTextPrices = [ordered]@{
  "D5" = '562.16'
  "D6" = '141.79'
  "D13" = '0.344'
  "D14" = '26.15'
  "D19" = '8.17'
  "D20" = '10.67'
  "D21" = '323.84'
  "D23" = '6.07'
  "D25" = '1.83'
  "D26" = '64.95'
  "D27" = '571.11'
  "D28" = '8.08'
  "D31" = '8.09'
  "D34" = '0.132'
  "D36" = '1.46'
  "D37" = '153.42'
  "D40" = '18.32'
  "D41" = '5.17'
  "D44" = '41.96'
  "D47" = '142.08'
  "D49" = '0.588'
  "D50" = '0.0507'
  "D51" = '19.34'
}
foreach ($addr in $forceTextPrices.Keys) {
  $cell = $ws.Range($addr)
  $cell.NumberFormat = "@"
  $cell.Value = $forceTextPrices[$addr]
  $cell.Style = "Normal"
}

# --- Price (column D): values that already contain more than one "." (or a
# subscript digit), so Excel stores them as text automatically.
$plainPrices = [ordered]@{
  "D2" = '60.874.44'
  "D3" = '2.402.24'
  "D9" = '2.410.15'
  "D17" = '60.730.51'
  "D18" = '2.401.84'
  "D30" = '0.0₃0941'
  "D46" = '0.0₆0279'
}
foreach ($addr in $plainPrices.Keys) {
  $ws.Range($addr).Value = $plainPrices[$addr]
}

# --- Volume(1h) (column E) for every changed row. Values are plain text
# padded with spaces, so Excel never treats them as numeric percentages.
$volumes = [ordered]@{
  "E2" = '  -0.30%  '
  "E3" = '  -1.05%  '
  "E4" = '  +0.44%  '
  "E5" = '  -1.84%  '
  "E6" = '  +1.12%  '
  "E7" = '  -0.31%  '
  "E8" = '  +1.46%  '
  "E9" = '  -0.18%  '
  "E10" = '  +0.17%  '
  "E11" = '  -0.45%  '
  "E12" = '  +1.35%  '
  "E13" = '  +1.27%  '
  "E14" = '  -0.35%  '
  "E15" = '  -1.39%  '
  "E17" = '  -0.33%  '
  "E18" = '  -0.27%  '
  "E19" = '  +6.57%  '
  "E20" = '  -0.29%  '
  "E21" = '  +0.01%  '
  "E22" = '  +0.61%  '
  "E23" = '  -0.57%  '
  "E24" = '  -0.23%  '
  "E25" = '  -2.99%  '
  "E26" = '  -0.10%  '
  "E27" = '  -2.75%  '
  "E28" = '  -4.93%  '
  "E30" = '  -0.09%  '
  "E31" = '  +2.00%  '
  "E32" = '  -2.18%  '
  "E33" = '  -2.60%  '
  "E34" = '  -0.01%  '
  "E35" = '  -0.58%  '
  "E36" = '  +2.96%  '
  "E37" = '  +0.93%  '
  "E38" = '  +0.32%  '
  "E39" = '  -1.37%  '
  "E40" = '  -0.12%  '
  "E41" = '  -0.13%  '
  "E44" = '  +1.82%  '
  "E45" = '  -0.67%  '
  "E46" = '  -3.85%  '
  "E47" = '  -0.35%  '
  "E48" = '  -0.08%  '
  "E49" = '  -0.73%  '
  "E50" = '  +0.18%  '
  "E51" = '  -2.04%  '
}
foreach ($addr in $volumes.Keys) {
  $ws.Range($addr).Value = $volumes[$addr]
}

# --- Rows 42/43: dogwifhat and USDe swapped rank positions; refresh their
# Coin name, Link, Price and Volume(1h) accordingly.
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.57'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +8.10%  '

$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.08%  '

Write-Host "Applied cryptos update."
